# LOM3049 - Termodinâmica de Máquinas: content refresh (build 2021-01-29)
#
# 1) Bump the "Ativação" date from 2018 to 2021.
# 2) Replace the "Objetivos" paragraph text.
# 3) Replace the "Programa resumido" paragraph text.
# 4) Replace the "Programa" paragraph text.

$d = $word.ActiveDocument

$wdReplaceAll = 2

function Replace-DocText($find, $replace) {
    $found = $d.Content.Find.Execute(
        $find, $true, $false, $false, $false, $false,
        $true, 1, $false, $replace, $wdReplaceAll)
    if (-not $found) {
        throw "Find.Execute could not locate text: $find"
    }
}

Replace-DocText `
    "Ativação: 01/01/2018" `
    "Ativação: 01/01/2021"

Replace-DocText `
    "1) Apresentar as 1ª e 2ª Leis da Termodinâmica e aplicá-las a problemas reais de engenharia;2) Calcular ciclos térmicos, ciclos de refrigeração e combustão, para que o Engenheiro de Materiais possa otimizar a eficiência de ciclos térmicos usando materiais que se adequem ás condições de projeto dos ciclos." `
    "Abordar os princípios básicos da termodinâmica de forma que os estudantes e futuros engenheiros tenham um entendimento claro e sólido sobre estes princípios. Apresentar diversos exemplos de engenharia do mundo real e de como a termodinâmica é aplicada na prática de engenharia. Enfatizar a compreensão da termodinâmica baseada na Física e em argumentos físicos, buscando incentivar o entendimento mais profundo da termodinâmica."

Replace-DocText `
    "1. Conceitos, Definições e Propriedades de uma substância pura2. Trabalho e Calor3. 1ª Lei de Termodinâmica4. 2ª Lei da Termodinâmica5. Entropia6. Ciclo Motores e de Refrigeração7. Projeto sobre Geração de Energia Termoelétrica: Ciclo Simples (vapor), Ciclo Combinado (turbina a gás/caldeira - turbina a vapor), Ciclos de refrigeração e de geração de potência combinados." `
    "1. Termodinâmica e Energia. 2. Importância das unidades e análise dimensional.3. Sistemas e volumes de controle. 4. Equipamentos domésticos e a Termodinâmica. 5. Propriedades de um sistema: estados termodinâmicos e equilíbrio. 6. Eficiência na conversão de energia. 7. Processos e ciclos térmicos. 8. Termodinâmica e o meio ambiente."

Replace-DocText `
    "1.Conceitos, definições e propriedades de uma substância pura; 2.Trabalho e calor;  3.1ª Lei da termodinâmica: Teoria e aplicação a volumes de controle; 4.2ª Lei da termodinâmica: Entropia5.2ª Lei da termodinâmica: Aplicação a volumes de controle;6.Ciclos motores Ciclos de refrigeração;7.Projeto sobre ciclo simples: Vapor; Projeto sobre ciclos combinados: Turbina a gás, turbina a vapor, Ciclos de refrigeração e de geração de potência combinados" `
    "1. Termodinâmica e Energia: formas de energia e transferência de energia por calor e trabalho; formas mecânicas de trabalho. 2. Sistema de Unidades e Análise Dimensional: importância na engenharia de máquinas. 3. Sistemas e volumes de controle: dispositivos ativos e passivos. 4. Propriedades de um sistema. Estados e equilíbrio: diagramas de propriedades para processos com mudança de fase; equilíbrio de estado do gás ideal; fator de compressibilidade; pressão de vapor e pressão de equilíbrio; calores específicos. 5. Balanço de energia em sistemas fechados e em volumes de controle: trabalho de fluxo e energia de escoamento de um fluido; regime permanente e transiente. 6. Máquinas térmicas e refrigeradores e a 2ª. Lei da Termodinâmica: princípios e ciclos de Carnot; entropia e variação de entropia em sólidos, líquidos e gases. 7. Eficiência na conversão de energia. Eficiência térmica. Eficiência de máquinas. Eficiência isoentrópica em dispositivos com escoamento em regime permanente. Balanço de entropia. 8. Processo e ciclos: Ciclos de potência a gás: Otto, Diesel, Stirling, Ericsson, Brayton e suas variações. Ciclos de potência a vapor e ciclos combinados gás-vapor: Rankine ideal; afastamento da condição ideal; eficiência do ciclo Rankine com e sem modificações; cogeração. Ciclos de refrigeração e sistemas de bombas de calor: sistemas a gás e por absorção. 9. Economia de energia: benefícios ao meio ambiente."

Write-Output "Done: 4 replacements applied."
